$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price/Volume columns to Text so numeric-looking strings are not
# auto-converted to numbers when the values are assigned below.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = '63.271.30'
$ws.Range("E2").Value = '  +0.22%  '
$ws.Range("D3").Value = '2.575.11'
$ws.Range("E3").Value = '  +0.37%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").Value = '588.91'
$ws.Range("E5").Value = '  +0.79%  '
$ws.Range("D6").Value = '144.67'
$ws.Range("E6").Value = '  -2.02%  '
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("D8").Value = '0.591'
$ws.Range("E8").Value = '  -1.45%  '
$ws.Range("D9").Value = '0.107'
$ws.Range("E9").Value = '  -1.92%  '
$ws.Range("D10").Value = '5.60'
$ws.Range("E10").Value = '  -0.73%  '
$ws.Range("D13").Value = '27.24'
$ws.Range("E13").Value = '  -0.59%  '
$ws.Range("D14").Value = '3.033.92'
$ws.Range("E14").Value = '  +0.28%  '
$ws.Range("D15").Value = '63.125.18'
$ws.Range("E15").Value = '  +0.07%  '
$ws.Range("E16").Value = '  -0.56%  '
$ws.Range("D17").Value = '2.569.87'
$ws.Range("E17").Value = '  +0.24%  '
$ws.Range("D18").Value = '11.06'
$ws.Range("E18").Value = '  -2.60%  '
$ws.Range("D19").Value = '340.55'
$ws.Range("E19").Value = '  -1.07%  '
$ws.Range("E20").Value = '  -1.81%  '
$ws.Range("E21").Value = '  -3.61%  '
$ws.Range("E22").Value = '  -0.02%  '
$ws.Range("E23").Value = '  +3.54%  '
$ws.Range("D24").Value = '67.87'
$ws.Range("E24").Value = '  +1.61%  '
$ws.Range("D25").Value = '1.57'
$ws.Range("E25").Value = '  +6.07%  '
$ws.Range("D26").Value = '1.63'
$ws.Range("E26").Value = '  -0.30%  '
$ws.Range("E27").Value = '  -2.98%  '
$ws.Range("E28").Value = '  -0.07%  '
$ws.Range("D29").Value = '7.94'
$ws.Range("E29").Value = '  -2.11%  '
$ws.Range("D30").Value = '8.24'
$ws.Range("E31").Value = '  -2.11%  '
$ws.Range("D32").Value = '472.62'
$ws.Range("E32").Value = '  +1.66%  '
$ws.Range("D33").Value = '0.0₃0803'
$ws.Range("E33").Value = '  -2.69%  '
$ws.Range("D34").Value = '1.69'
$ws.Range("E34").Value = '  +3.71%  '
$ws.Range("D35").Value = '176.28'
$ws.Range("E35").Value = '  +0.14%  '
$ws.Range("E36").Value = '  +0.05%  '
$ws.Range("E37").Value = '  -1.69%  '
$ws.Range("E38").Value = '  -1.95%  '
$ws.Range("D39").Value = '4.55'
$ws.Range("E39").Value = '  +0.08%  '
$ws.Range("E40").Value = '  -0.02%  '
$ws.Range("E41").Value = '  -2.90%  '
$ws.Range("D42").Value = '40.12'
$ws.Range("E42").Value = '  +1.25%  '
$ws.Range("D43").Value = '157.71'
$ws.Range("E43").Value = '  +4.29%  '
$ws.Range("D44").Value = '3.69'
$ws.Range("E44").Value = '  -3.51%  '
$ws.Range("D45").Value = '21.34'
$ws.Range("D46").Value = '0.634'
$ws.Range("E46").Value = '  +3.35%  '
$ws.Range("E47").Value = '  -1.60%  '
$ws.Range("D48").Value = '0.0963'
$ws.Range("E48").Value = '  -1.40%  '
$ws.Range("E49").Value = '  -1.19%  '
$ws.Range("D50").Value = '18.13'
$ws.Range("E50").Value = '  -1.78%  '
$ws.Range("E51").Value = '  -0.06%  '

# Restore original (unformatted) cell style now that the text values are set.
$dataRange.ClearFormats()

